# Update the 13x13 correlation matrix in corrM (Sheet1, range B2:N14)
# with the refreshed correlation values ("updated data and results till 2020 Jan.").
# The matrix is symmetric with 1s on the diagonal; row/column headers (row 1, col A)
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the new 13x13 matrix of correlation values (rows/cols correspond to B2:N14)
$newCorr = New-Object 'object[,]' 13,13

# row 2
$newCorr[0,0] = 1
$newCorr[0,1] = -0.680182012841634
$newCorr[0,2] = -0.0783303930584449
$newCorr[0,3] = -0.03052985330403633
$newCorr[0,4] = 0.01371389403032016
$newCorr[0,5] = 0.0124000437774399
$newCorr[0,6] = -0.04819288080125357
$newCorr[0,7] = -0.1023235157253121
$newCorr[0,8] = -0.01461638364964552
$newCorr[0,9] = 0.05652648764826258
$newCorr[0,10] = 0.02823333018647582
$newCorr[0,11] = -0.03196208928288601
$newCorr[0,12] = -0.01748914489945714

# row 3
$newCorr[1,0] = -0.680182012841634
$newCorr[1,1] = 1
$newCorr[1,2] = 0.1221720518472285
$newCorr[1,3] = 0.05653006462466296
$newCorr[1,4] = -0.05391471495687534
$newCorr[1,5] = -0.06253280915464397
$newCorr[1,6] = 0.1111497886989482
$newCorr[1,7] = 0.06909878822468789
$newCorr[1,8] = -0.01909015669427782
$newCorr[1,9] = -0.04374136016896989
$newCorr[1,10] = 0.006575284563791996
$newCorr[1,11] = 0.07679382384287423
$newCorr[1,12] = 0.03795934618164015

# row 4
$newCorr[2,0] = -0.0783303930584449
$newCorr[2,1] = 0.1221720518472285
$newCorr[2,2] = 1
$newCorr[2,3] = 0.2981098520604903
$newCorr[2,4] = -0.1633195338039325
$newCorr[2,5] = -0.1470984054862519
$newCorr[2,6] = 0.3588750911617425
$newCorr[2,7] = -0.5151684230962238
$newCorr[2,8] = 0.01043230032222128
$newCorr[2,9] = -0.459597172319483
$newCorr[2,10] = -0.342156576835674
$newCorr[2,11] = 0.3420214029909316
$newCorr[2,12] = -0.4634142680280373

# row 5
$newCorr[3,0] = -0.03052985330403633
$newCorr[3,1] = 0.05653006462466296
$newCorr[3,2] = 0.2981098520604903
$newCorr[3,3] = 1
$newCorr[3,4] = 0.3612509639485904
$newCorr[3,5] = 0.3641198327393359
$newCorr[3,6] = 0.4127631093249592
$newCorr[3,7] = 0.03599422413450309
$newCorr[3,8] = 0.7921024341773386
$newCorr[3,9] = 0.08673123342589122
$newCorr[3,10] = -0.002601885874268448
$newCorr[3,11] = 0.3904627113860717
$newCorr[3,12] = -0.0004828960342666396

# row 6
$newCorr[4,0] = 0.01371389403032016
$newCorr[4,1] = -0.05391471495687534
$newCorr[4,2] = -0.1633195338039325
$newCorr[4,3] = 0.3612509639485904
$newCorr[4,4] = 1
$newCorr[4,5] = 0.9809329076610995
$newCorr[4,6] = 0.3722554629941501
$newCorr[4,7] = 0.1806114662513572
$newCorr[4,8] = 0.6301873021039683
$newCorr[4,9] = -0.04454127211874876
$newCorr[4,10] = -0.3129315469833898
$newCorr[4,11] = 0.3725002900626034
$newCorr[4,12] = 0.05618209907557382

# row 7
$newCorr[5,0] = 0.0124000437774399
$newCorr[5,1] = -0.06253280915464397
$newCorr[5,2] = -0.1470984054862519
$newCorr[5,3] = 0.3641198327393359
$newCorr[5,4] = 0.9809329076610995
$newCorr[5,5] = 1
$newCorr[5,6] = 0.3743625636196435
$newCorr[5,7] = 0.2414502345040999
$newCorr[5,8] = 0.6426601479267363
$newCorr[5,9] = 0.008927636603904878
$newCorr[5,10] = -0.268192419546511
$newCorr[5,11] = 0.3823050395265182
$newCorr[5,12] = 0.07325264309598917

# row 8
$newCorr[6,0] = -0.04819288080125357
$newCorr[6,1] = 0.1111497886989482
$newCorr[6,2] = 0.3588750911617425
$newCorr[6,3] = 0.4127631093249592
$newCorr[6,4] = 0.3722554629941501
$newCorr[6,5] = 0.3743625636196435
$newCorr[6,6] = 1
$newCorr[6,7] = -0.1997910261620322
$newCorr[6,8] = 0.4475614826426157
$newCorr[6,9] = -0.1437392066062044
$newCorr[6,10] = -0.2180916373139349
$newCorr[6,11] = 0.9172288078615954
$newCorr[6,12] = -0.159014452899539

# row 9
$newCorr[7,0] = -0.1023235157253121
$newCorr[7,1] = 0.06909878822468789
$newCorr[7,2] = -0.5151684230962238
$newCorr[7,3] = 0.03599422413450309
$newCorr[7,4] = 0.1806114662513572
$newCorr[7,5] = 0.2414502345040999
$newCorr[7,6] = -0.1997910261620322
$newCorr[7,7] = 1
$newCorr[7,8] = 0.2333820891919775
$newCorr[7,9] = 0.6883971392785034
$newCorr[7,10] = 0.5246616647241267
$newCorr[7,11] = -0.2193972891311861
$newCorr[7,12] = 0.7480185820308105

# row 10
$newCorr[8,0] = -0.01461638364964552
$newCorr[8,1] = -0.01909015669427782
$newCorr[8,2] = 0.01043230032222128
$newCorr[8,3] = 0.7921024341773386
$newCorr[8,4] = 0.6301873021039683
$newCorr[8,5] = 0.6426601479267363
$newCorr[8,6] = 0.4475614826426157
$newCorr[8,7] = 0.2333820891919775
$newCorr[8,8] = 1
$newCorr[8,9] = 0.2711736763683567
$newCorr[8,10] = 0.05649703490390926
$newCorr[8,11] = 0.4819490113416408
$newCorr[8,12] = 0.1951778656562174

# row 11
$newCorr[9,0] = 0.05652648764826258
$newCorr[9,1] = -0.04374136016896989
$newCorr[9,2] = -0.459597172319483
$newCorr[9,3] = 0.08673123342589122
$newCorr[9,4] = -0.04454127211874876
$newCorr[9,5] = 0.008927636603904878
$newCorr[9,6] = -0.1437392066062044
$newCorr[9,7] = 0.6883971392785034
$newCorr[9,8] = 0.2711736763683567
$newCorr[9,9] = 1
$newCorr[9,10] = 0.931601405267243
$newCorr[9,11] = -0.13932507739623
$newCorr[9,12] = 0.759344995191208

# row 12
$newCorr[10,0] = 0.02823333018647582
$newCorr[10,1] = 0.006575284563791996
$newCorr[10,2] = -0.342156576835674
$newCorr[10,3] = -0.002601885874268448
$newCorr[10,4] = -0.3129315469833898
$newCorr[10,5] = -0.268192419546511
$newCorr[10,6] = -0.2180916373139349
$newCorr[10,7] = 0.5246616647241267
$newCorr[10,8] = 0.05649703490390926
$newCorr[10,9] = 0.931601405267243
$newCorr[10,10] = 1
$newCorr[10,11] = -0.2257734971776478
$newCorr[10,12] = 0.7016434961541568

# row 13
$newCorr[11,0] = -0.03196208928288601
$newCorr[11,1] = 0.07679382384287423
$newCorr[11,2] = 0.3420214029909316
$newCorr[11,3] = 0.3904627113860717
$newCorr[11,4] = 0.3725002900626034
$newCorr[11,5] = 0.3823050395265182
$newCorr[11,6] = 0.9172288078615954
$newCorr[11,7] = -0.2193972891311861
$newCorr[11,8] = 0.4819490113416408
$newCorr[11,9] = -0.13932507739623
$newCorr[11,10] = -0.2257734971776478
$newCorr[11,11] = 1
$newCorr[11,12] = -0.1896858643044187

# row 14
$newCorr[12,0] = -0.01748914489945714
$newCorr[12,1] = 0.03795934618164015
$newCorr[12,2] = -0.4634142680280373
$newCorr[12,3] = -0.0004828960342666396
$newCorr[12,4] = 0.05618209907557382
$newCorr[12,5] = 0.07325264309598917
$newCorr[12,6] = -0.159014452899539
$newCorr[12,7] = 0.7480185820308105
$newCorr[12,8] = 0.1951778656562174
$newCorr[12,9] = 0.759344995191208
$newCorr[12,10] = 0.7016434961541568
$newCorr[12,11] = -0.1896858643044187
$newCorr[12,12] = 1

# Write the whole matrix back into the worksheet in one shot
$ws.Range("B2:N14").Value = $newCorr

Write-Host "Updated corrM matrix (B2:N14) with refreshed correlation values."
